# Implicit evaluations study codebook - AMP and BFI, IAT Task completed
# Tidy up the wording of the explanation column (B) for the participant-level
# variables: capitalize the first letter and end each sentence with a period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A code so that the various data can be assigned to the respective participants - every participant has an own code."
$ws.Range("B3").Value = "The date on which the data was collected from the respective participant."
$ws.Range("B4").Value = "The exact time when the participants took part in the experiment. "
$ws.Range("B5").Value = "The age of the respective participant (given in years!)."
$ws.Range("B6").Value = "The gender that the respective participant identifies as."

# Leave the cursor where the author left it when they saved the file.
$ws.Range("I7").Select() | Out-Null
